$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: new data row - "م" (id), name, current-balance ratio, price, transaction-count ratio
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "HELI-CURE 14 ENTERIC COATED TAB"
$ws.Range("H4").Value = "2:1"
$ws.Range("L4").Value = 120
$ws.Range("N4").Value = "0:2"

# Row 5: totals row gains a value and grows slightly taller
$ws.Range("K5").Value = 120
$ws.Rows(5).RowHeight = 26.25
